$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column Z (Augment Size) for all data rows from numeric 20 to text "-"
$ws.Range("Z2:Z257").Value = "-"

# Update the active selection to match the saved view state
$ws.Activate()
$ws.Range("M230").Select()
